$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Rows 29-36: shift the September (R/S) entries down by one,
# inserting a brand-new "internet verify it" entry at the top (row 29).
$ws.Range("R29").Value = "internet verify it"
$ws.Range("S29").Value = "2024-09-03 19:05:49"

$ws.Range("R30").Value = "balance your axis"
$ws.Range("S30").Value = "2024-09-03 13:14:06"

$ws.Range("R31").Value = "lounge"
$ws.Range("S31").Value = "2024-09-03 13:08:08"

$ws.Range("R32").Value = "balance your axis"
$ws.Range("S32").Value = "2024-09-03 11:21:30"

$ws.Range("R33").Value = "broker"
$ws.Range("S33").Value = "2024-09-01 22:35:38"

$ws.Range("R34").Value = "amazeloan"
$ws.Range("S34").Value = "2024-09-01 10:12:03"

$ws.Range("R35").Value = "amazeloan"
$ws.Range("S35").Value = "2024-09-01 09:42:38"

$ws.Range("R36").Value = "amazeloan"
$ws.Range("S36").Value = "2024-09-01 09:29:24"

# Row 37 previously held an August (P/Q) "hdfc" entry; that entry shifts
# down into row 38, and row 37 becomes the last amazeloan (September) row.
$ws.Range("P37").Value = ""
$ws.Range("Q37").Value = ""
$ws.Range("R37").Value = "amazeloan"
$ws.Range("S37").Value = "2024-09-01 09:27:06"

# Rows 38-40: shift the August (P/Q) "hdfc" entries down by one.
$ws.Range("P38").Value = "hdfc"
$ws.Range("Q38").Value = "2024-08-30 12:15:48"

$ws.Range("P39").Value = "hdfc"
$ws.Range("Q39").Value = "2024-08-21 20:17:10"

$ws.Range("P40").Value = "hdfc"
$ws.Range("Q40").Value = "2024-08-21 20:16:45"

# Row 41 used to hold the "Broadband" group label; it now holds the last
# shifted hdfc entry, and "Broadband" moves down to the new row 42.
$ws.Range("A41").Value = ""
$ws.Range("P41").Value = "hdfc"
$ws.Range("Q41").Value = "2024-08-21 20:15:50"

$ws.Range("A42").Value = "Broadband"
$ws.Range("B42:Y42").Value = ""

